$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.202.22"
$ws.Range("E2").Value = "  +9.29%  "
$ws.Range("D3").Value = "2.627.82"
$ws.Range("E3").Value = "  +12.29%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "510.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.46"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +9.37%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.32%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.610"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.34%  "
$ws.Range("D9").Value = "2.627.18"
$ws.Range("E9").Value = "  +11.86%  "
$ws.Range("E10").Value = "  +12.72%  "
$ws.Range("E11").Value = "  +8.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.347"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.07%  "
$ws.Range("E13").Value = "  +1.73%  "
$ws.Range("D14").Value = "3.056.02"
$ws.Range("E14").Value = "  +11.10%  "
$ws.Range("D15").Value = "60.119.57"
$ws.Range("E15").Value = "  +9.10%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "22.18"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +11.09%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000140"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +8.80%  "
$ws.Range("D18").Value = "2.611.34"
$ws.Range("E18").Value = "  +11.29%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +10.06%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +10.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +9.34%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +7.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.424"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +7.72%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.169"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +12.60%  "
$ws.Range("D27").Value = "2.685.98"
$ws.Range("E27").Value = "  +9.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.990"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.87%  "
$ws.Range("D29").Value = "0.0₃0855"
$ws.Range("E29").Value = "  +15.92%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "158.19"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.66"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +9.17%  "
$ws.Range("E34").Value = "  +8.13%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.59"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +10.20%  "
$ws.Range("E36").Value = "  +11.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.97"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +10.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "314.39"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +25.96%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.874"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.14%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.48"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +11.91%  "
$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.64%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "35.45"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.637"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.54%  "
$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0578"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +12.05%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.102"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.792"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +28.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.995"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +19.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.01"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +14.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0238"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +8.42%  "
$ws.Range("D51").Value = "2.005.45"
$ws.Range("E51").Value = "  +12.08%  "
